$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 104

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null) {
        $s = $val.ToString()
        if ($s.EndsWith("16")) {
            $cell.Value2 = $s.Substring(0, $s.Length - 2)
        }
    }
}
